$wb = $excel.ActiveWorkbook

# --- mon (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = "CSC221"
$ws.Range("D2").Value = "CSC221"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("J4").Value = "GST121"
$ws.Range("K4").Value = "GST121"
$ws.Range("F5").Value = "MIS421"
$ws.Range("G5").Value = "MIS421"
$ws.Range("K6").Value = "MAT226"
$ws.Range("E7").Value = ""
$ws.Range("E11").Value = "CIT242"
$ws.Range("F11").Value = "CIT242"
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("H13").Value = "CSC227"
$ws.Range("I13").Value = "CSC227"
$ws.Range("J16").Value = "PHY122"
$ws.Range("K16").Value = "PHY122"
$ws.Range("E17").Value = "CSC121"
$ws.Range("F17").Value = "CSC121"
$ws.Range("H18").Value = "MIS221"
$ws.Range("H21").Value = "CSC223"
$ws.Range("J24").Value = ""
$ws.Range("K24").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("K25").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("J27").Value = "CIS421"
$ws.Range("K27").Value = "CIS421"
$ws.Range("I29").Value = "CSC442"
$ws.Range("J29").Value = "CSC442"

# --- tue (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("E2").Value = "EDS121"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("J4").Value = "MAT121"
$ws.Range("K4").Value = "MAT121"
$ws.Range("D5").Value = "MAT225"
$ws.Range("E5").Value = "MAT225"
$ws.Range("H6").Value = "CSC446"
$ws.Range("I6").Value = "CSC446"
$ws.Range("B7").Value = "CSC444"
$ws.Range("C7").Value = "CSC444"
$ws.Range("G7").Value = "DLD121"
$ws.Range("H7").Value = "DLD121"
$ws.Range("J8").Value = "CSC225"
$ws.Range("K8").Value = "CSC225"
$ws.Range("I9").Value = "BUS124"
$ws.Range("J9").Value = "BUS124"
$ws.Range("F11").Value = "MIS425"
$ws.Range("G11").Value = "MIS425"
$ws.Range("G13").Value = "MAT229"
$ws.Range("H13").Value = "MAT229"
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = "PHY121"
$ws.Range("K14").Value = "PHY121"
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("J18").Value = "CBS121"
$ws.Range("K18").Value = ""
$ws.Range("C19").Value = "MIS426"
$ws.Range("D19").Value = "MIS426"
$ws.Range("D20").Value = "CSC226"
$ws.Range("E20").Value = "CSC226"
$ws.Range("H20").Value = "CSC424"
$ws.Range("I20").Value = "CSC424"
$ws.Range("I21").Value = "CSC423"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("H26").Value = ""
$ws.Range("I26").Value = "CSC224"
$ws.Range("J26").Value = "CSC224"
$ws.Range("F29").Value = ""
$ws.Range("J29").Value = "MIS423"
$ws.Range("K29").Value = "MIS423"

# --- wed (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("G2").Value = "CIT141"
$ws.Range("H2").Value = "CIT141"
$ws.Range("B3").Value = ""
$ws.Range("D4").Value = "TMC121"
$ws.Range("C6").Value = "CSC223"
$ws.Range("D6").Value = "CSC223"
$ws.Range("I7").Value = "DLD221"
$ws.Range("J7").Value = "DLD221"
$ws.Range("K7").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("C13").Value = "EDS421"
$ws.Range("F14").Value = "CSC425"
$ws.Range("H14").Value = "CSC441"
$ws.Range("I14").Value = "CSC441"
$ws.Range("C15").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = "PHY129"
$ws.Range("I25").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("H26").Value = ""
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = ""
$ws.Range("B28").Value = "BUS326"
$ws.Range("C28").Value = "BUS326"
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = ""

# --- thur (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("I2").Value = "CIT221"
$ws.Range("J2").Value = "CIT221"
$ws.Range("J4").Value = "CSC125"
$ws.Range("K4").Value = "CSC125"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I11").Value = "MIS421"
$ws.Range("F13").Value = "CSC423"
$ws.Range("G13").Value = "CSC423"
$ws.Range("D15").Value = "ACC121"
$ws.Range("E15").Value = "ACC121"
$ws.Range("H15").Value = ""
$ws.Range("G17").Value = "CIT121"
$ws.Range("H17").Value = "CIT121"
$ws.Range("H18").Value = "MIS423"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = "MAT122"
$ws.Range("E20").Value = "MAT122"
$ws.Range("F21").Value = "CSC425"
$ws.Range("G21").Value = "CSC425"
$ws.Range("I21").Value = "EDS221"
$ws.Range("F23").Value = "BUS221"
$ws.Range("G23").Value = "BUS221"
$ws.Range("H24").Value = ""
$ws.Range("I24").Value = ""
$ws.Range("I27").Value = ""
$ws.Range("J27").Value = ""

# --- fri (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("E2").Value = "CSC225"
$ws.Range("C3").Value = "MAT226"
$ws.Range("D3").Value = "MAT226"
$ws.Range("G8").Value = "MAT225"
$ws.Range("E12").Value = "BUS326"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("E13").Value = "CIT224"
$ws.Range("F13").Value = "CIT224"
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("F20").Value = "GST222"
$ws.Range("G20").Value = "GST222"
$ws.Range("C21").Value = "MAT121"
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = ""
$ws.Range("B28").Value = "CBS221"
$ws.Range("E29").Value = ""
